# Automated daily update of the tracking sheet.
# For each data row (columns: D=total days, E=days remaining, F=start date as yyyymmdd):
#   - compute the cycle's end date = F + D - 1 days
#   - if "today" (2025-11-17) is past the end date, the cycle has expired:
#       reset E back to the full total (D) and set F to today's date
#   - otherwise the cycle is still running: just decrement E by 1 (one day elapsed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = Get-Date -Year 2025 -Month 11 -Day 17

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    $fStr = [string][int]$fVal
    if ($fStr.Length -ne 8) {
        # malformed/unexpected date value, skip this row untouched
        continue
    }

    $fYear = [int]$fStr.Substring(0, 4)
    $fMonth = [int]$fStr.Substring(4, 2)
    $fDay = [int]$fStr.Substring(6, 2)

    $startDate = Get-Date -Year $fYear -Month $fMonth -Day $fDay
    $endDate = $startDate.AddDays([int]$dVal - 1)

    if ($endDate -lt $today) {
        $newE = [int]$dVal
        $newFDate = $today
    } else {
        $newE = [int]$eVal - 1
        $newFDate = $startDate
    }

    $newF = [int]$newFDate.ToString("yyyyMMdd")

    $eCell.Value = $newE
    $fCell.Value = $newF
}
